$wb = $excel.ActiveWorkbook

# --- Update view state (selection) on existing sheets ---

$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Activate()
$wsGermany.Range("A3").Select()

$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsBelgium.Activate()
$wsBelgium.Range("A8:A35").Select()

$wsCzech = $wb.Worksheets.Item("Czech")
$wsCzech.Activate()
$wsCzech.Range("A25:A31").Select()

$wsPortugal = $wb.Worksheets.Item("Portugal")
$wsPortugal.Activate()
$wsPortugal.Range("A32").Select()

$wsSlovakia = $wb.Worksheets.Item("Slovakia")
$wsSlovakia.Activate()
$wsSlovakia.Range("B18").Select()

# --- Add the new "Italy" sheet, based on the "Czech" sheet layout ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCzech.Copy($null, $lastSheet)
$wsItaly = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItaly.Name = "Italy"

$wsItaly.Range("B2").Value = "Italy Market"
$wsItaly.Range("B4").Value = "NGC-3145/T2153/T2218/T2456"
$wsItaly.Range("B4").Style = "Normal"
$wsItaly.Range("A32").Value = "MX4000"

$wsItaly.Activate()
$wsItaly.Range("B4").Select()
